$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admission_High School Required")
$ws.Rows.Item(12).Delete()
